# Revert "Fix max righe/colonne - other"
#
# The original commit (now being reverted) had:
#  - filled in a "Done/Date/Note" row for the Question/Answer items (row 3 & 4)
#  - added a "Mancante: URL" note in E4
#  - resized several rows (presumably as part of a "max rows/cols" layout fix)
#
# This script undoes those edits: clears the newly-added cell values (which in
# turn drops the now-unused "Mancante: URL" shared string on save) and restores
# the previous row heights / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the cell content that was added by the commit being reverted ---
# Row 3: "Done"/"Data" values for the Question/Answer task
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
# Row 4: "Done"/"Data" values + the "Mancante: URL" note
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# --- Restore the previous row heights ---
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 144
$ws.Rows.Item(7).RowHeight = 28.8
# Rows 11 & 12 go back to the sheet's default (no explicit custom height)
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(17).RowHeight = 28.8
$ws.Rows.Item(18).RowHeight = 28.8
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 43.2
$ws.Rows.Item(25).RowHeight = 28.8
$ws.Rows.Item(26).RowHeight = 28.8
$ws.Rows.Item(27).RowHeight = 115.2

# --- Restore the previous selection ---
$ws.Range("B6").Select()
